$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "JPH7"
$ws.Range("B12").Value = "app"
$ws.Range("A13").Value = "JP7"
$ws.Range("B13").Value = "grant"

$ws.Range("B13").Select()
